# Delete the empty placeholder sheet "Лист1" and rename the data sheet
# "Sheet1" to "Result".

$wb = $excel.ActiveWorkbook

$excel.DisplayAlerts = $false

$emptySheet = $wb.Worksheets.Item("Лист1")
$emptySheet.Delete()

$dataSheet = $wb.Worksheets.Item("Sheet1")
$dataSheet.Name = "Result"

$dataSheet.Range("D2:E2").Copy()
$dataSheet.Range("D3:E8").PasteSpecial(-4122)  # xlPasteFormats

# Slightly adjust column widths (A narrower, C/D/E wider)
$dataSheet.Columns.Item(1).ColumnWidth = 19.166666666666668
$dataSheet.Columns.Item(2).ColumnWidth = 15.333333333333334
$dataSheet.Columns.Item(3).ColumnWidth = 14.0
$dataSheet.Columns.Item(4).ColumnWidth = 15.333333333333334
$dataSheet.Columns.Item(5).ColumnWidth = 17.333333333333336
$dataSheet.Columns.Item(6).ColumnWidth = 18.5

$dataSheet.Activate()
$dataSheet.Range("D12").Select()

$excel.DisplayAlerts = $true
